$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the authoritative diff.
# Values that look like plain numbers are prefixed with a literal single quote
# so Excel stores them as text (matching the original inline-string cells)
# instead of auto-converting them to numeric cells.
$updates = @{
    "D2" = '22.939.48'
    "E2" = '  -1.11%  '
    "D3" = '1.574.32'
    "E3" = '  -2.09%  '
    "D4" = "'1.003"
    "E4" = '  +0.28%  '
    "E5" = '  +0.12%  '
    "D6" = "'299.22"
    "E6" = '  -1.25%  '
    "D7" = "'0.3739"
    "E7" = '  -0.74%  '
    "D8" = "'0.3545"
    "E8" = '  -2.90%  '
    "D9" = "'49.87"
    "E9" = '  +2.52%  '
    "D10" = "'1.004"
    "E10" = '  +0.29%  '
    "D11" = "'1.212"
    "E11" = '  -4.44%  '
    "D12" = "'0.07953"
    "E12" = '  -1.60%  '
    "D13" = "'21.74"
    "E13" = '  -5.56%  '
    "D14" = "'6.406"
    "E14" = '  -2.67%  '
    "D15" = "'7.271"
    "E15" = '  -4.81%  '
    "D16" = "'0.00001220"
    "E16" = '  -3.86%  '
    "D17" = '1.575.57'
    "E17" = '  -2.11%  '
    "D18" = "'91.71"
    "D19" = "'0.06732"
    "E19" = '  -0.80%  '
    "D20" = "'17.62"
    "E20" = '  -4.00%  '
    "E21" = '  +0.13%  '
    "D22" = "'6.356"
    "E22" = '  -3.40%  '
    "D23" = '22.944.95'
    "E23" = '  -1.14%  '
    "D24" = "'12.57"
    "E24" = '  -4.17%  '
    "D25" = "'2.370"
    "E25" = '  +0.74%  '
    "D26" = "'2.796"
    "E26" = '  -3.83%  '
    "D27" = "'20.54"
    "D28" = "'147.29"
    "E28" = '  -2.02%  '
    "D29" = "'5.151"
    "E29" = '  -2.21%  '
    "D30" = "'131.33"
    "E30" = '  -0.85%  '
    "D31" = "'2.327"
    "E31" = '  -3.80%  '
    "D32" = "'6.511"
    "E32" = '  -5.96%  '
    "D33" = '1.751.78'
    "E33" = '  -1.89%  '
    "D34" = "'0.9303"
    "E34" = '  -4.88%  '
    "D35" = "'0.07336"
    "E35" = '  -5.08%  '
    "D36" = "'0.08747"
    "E36" = '  -1.26%  '
    "B37" = 'FraxShare'
    "C37" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "D37" = "'9.916"
    "E37" = '  -1.79%  '
    "B38" = 'VeChain'
    "C38" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "D38" = "'0.02627"
    "E38" = '  -5.68%  '
    "D39" = "'0.2458"
    "E39" = '  -3.64%  '
    "D40" = "'5.962"
    "E40" = '  -4.89%  '
    "D41" = "'1.339"
    "E41" = '  -4.06%  '
    "D42" = "'0.6844"
    "E42" = '  -4.40%  '
    "D43" = "'11.81"
    "E43" = '  -7.62%  '
    "D44" = "'14.73"
    "E44" = '  -7.13%  '
    "D45" = "'1.000"
    "E45" = '  +0.02%  '
    "D46" = "'0.6313"
    "E46" = '  -4.41%  '
    "D47" = "'3.964"
    "E47" = '  -0.49%  '
    "D48" = "'2.239"
    "E48" = '  -2.63%  '
    "D49" = "'130.07"
    "E49" = '  -0.96%  '
    "D50" = "'0.07845"
    "E50" = '  -2.09%  '
    "D51" = "'1.180"
    "E51" = '  +1.01%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
